$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.9552290439605713, 84.02)
    3  = @(0.2141249179840088, 38.991)
    4  = @(1.249463796615601, 124.622)
    5  = @(0.4011049270629883, 88.65900000000001)
    6  = @(0.3084449768066406, 6.017)
    7  = @(0.1787829399108887, 6.017)
    8  = @(0.2351679801940918, 6.059)
    9  = @(0.2315452098846436, 6.059)
    10 = @(0.1913321018218994, 7.972)
    11 = @(0.1875576972961426, 11.609)
    12 = @(0.2387080192565918, 10.307)
    13 = @(0.2389688491821289, 8.481999999999999)
    14 = @(0.1888337135314941, 201.897)
    15 = @(0.205280065536499, 522.954)
    16 = @(0.5307791233062744, 65.334)
    17 = @(0.6413018703460693, 528.609)
    18 = @(0.1999771595001221, 6.215)
    19 = @(0.2328498363494873, 6.215)
    20 = @(0.2479281425476074, 6.077)
    21 = @(0.2545340061187744, 6.077)
    22 = @(0.2219791412353516, 8.17)
    23 = @(0.301379919052124, 8.17)
    24 = @(0.2526941299438477, 8.032)
    25 = @(0.2814657688140869, 8.032)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}
